# Journal de travail - add two new entries (rows 25 and 26) and move selection
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: 11/03/2018 - Rédaction des conventions (0.5h)
$ws.Range("A25").Value = 43170
$ws.Range("B25").Value = "Rédaction des conventions à employer pour le code et les commentaires Java"
$ws.Range("C25").Value = 0.5

# Row 26: 12/03/2018 - Modélisation du schéma relationnel (0.75h)
$ws.Range("A26").Value = 43171
$ws.Range("B26").Value = "Modélisation du schéma relationnel sur papier avec Héléna"
$ws.Range("C26").Value = 0.75

# Row 26 needs a taller row to fit the wrapped text, matching the other
# two-line entries in the sheet.
$ws.Rows.Item(26).RowHeight = 30

# Move the active selection to the next empty row, as in the saved file.
$ws.Range("A28").Select() | Out-Null
